# Apply committee-sheet updates:
#  * Fix affiliation typo for Agnieszka Kubik-Komar
#  * Insert new committee member (Lilla Di Scala, Johnson & Johnson) as row 22,
#    pushing the existing "local" committee rows down by one

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct "University of Life Sciences in Lubli" -> "...Lublin"
$ws.Range("C18").Value = "University of Life Sciences in Lublin"

# Insert a new row above the current row 22 (Laurence Giullier) and fill it
# in with the new committee member
$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = "Lilla"
$ws.Range("B22").Value = "Di Scala"
$ws.Range("C22").Value = "Johnson & Johnson"

# Leave the selection where the author's session ended up
$ws.Range("I14").Select()
